$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "pinna di squalo" (shark fin) saturation data: new rows 14-22, cols M:P ---
$data = @{
    14 = @(699, 6, 3546, 30)
    15 = @(799, 6, 3819, 31)
    16 = @(899, 7, 3845, 32)
    17 = @(999, 8, 3860, 32)
    18 = @(1199, 9, 3855, 32)
    19 = @(1399, 11, 3883, 32)
    20 = @(1599, 12, 3890, 32)
    21 = @(1799, 14, 3906, 32)
    22 = @(1999, 15, 3921, 32)
}

foreach ($row in 14..22) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 13).Value = $vals[0]   # M
    $ws.Cells.Item($row, 14).Value = $vals[1]   # N
    $ws.Cells.Item($row, 15).Value = $vals[2]   # O
    $ws.Cells.Item($row, 16).Value = $vals[3]   # P
    $ws.Cells.Item($row, 17).Formula = "=O$row/M$row"
    $ws.Cells.Item($row, 18).Formula = "=Q$row*SQRT((N$row/M$row)^2 + (P$row/O$row)^2)"
}

# R21 (saturated point, "pinna di squalo") gets the underline style, like B15
$ws.Range("R21").Font.Underline = $true

# Summary stats over the non-saturated shark-fin region O16:O22
$ws.Range("O23").Formula = "=AVERAGE(O16:O22)"
$ws.Range("P23").Formula = "=STDEV(O16:O22)"
$ws.Range("O23").Font.Bold = $true

# Reposition the view: scrolled down with R21 selected
$ws.Range("R21").Select()

Write-Output "edit applied"
